$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1. The "_GoBack" bookmark used to sit alone in the trailing empty
#    paragraph at the end of the document. Word moves this bookmark to
#    track the last edit location, so we delete it from its old spot
#    first (it gets re-created at the new edit point below).
# ------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# ------------------------------------------------------------------
# 2. Split the "Bensound.com" hyperlink run into "Bensound" + ".com"
#    runs (identical run formatting on both sides) with the _GoBack
#    bookmark now placed at the split point, still inside the
#    <w:hyperlink> that points at https://www.bensound.com/ (rId4).
# ------------------------------------------------------------------
$hyperlink = $d.Hyperlinks.Item(1)
$linkRange = $d.Range($hyperlink.Range.Start, $hyperlink.Range.End)

$splitXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' + `
'<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' + `
'<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData>' + `
'<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:r="http://schemas.openxmlformats.org/officeDocument/2006/relationships">' + `
'<w:body><w:p>' + `
'<w:hyperlink r:id="rId4" w:history="1">' + `
'<w:r><w:rPr><w:rStyle w:val="Hyperlink"/><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:b/><w:bCs/><w:color w:val="9CD121"/><w:bdr w:val="none" w:sz="0" w:space="0" w:color="auto" w:frame="1"/></w:rPr><w:t>Bensound</w:t></w:r>' + `
'<w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/>' + `
'<w:r><w:rPr><w:rStyle w:val="Hyperlink"/><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:b/><w:bCs/><w:color w:val="9CD121"/><w:bdr w:val="none" w:sz="0" w:space="0" w:color="auto" w:frame="1"/></w:rPr><w:t>.com</w:t></w:r>' + `
'</w:hyperlink>' + `
'</w:p></w:body></w:document>' + `
'</pkg:xmlData></pkg:part></pkg:package>'

$linkRange.InsertXML($splitXml)

# ------------------------------------------------------------------
# 3. Add the (now used) FollowedHyperlink character style to
#    styles.xml, matching Word's built-in definition as closely as
#    the object model allows.
# ------------------------------------------------------------------
$followed = $d.Styles.Add("FollowedHyperlink", 2)
$followed.BaseStyle = $d.Styles("DefaultParagraphFont")
$followed.Priority = 99
$followed.UnhideWhenUsed = $true
$followed.Font.Color = 7491477
$followed.Font.Underline = 1
